# Split March expenses onto a sheet named after its month ("mars"), and
# start tracking April's expenses on a new sheet ("april").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "mars"

# Remove the bold header styling that used to mark row 1.
$ws.Range("A1:E1").ClearFormats()

# Append the new March rows that were entered after the original edit.
$marsRows = @(
    @("Food", "sws", "2023-03-15", "123.0", "Savings"),
    @("Transportation", "sws", "2023-03-15", "343.0", "Savings"),
    @("Transportation", "sws", "2024-03-13", "5454.0", "Checkings")
)

$r = 16
foreach ($row in $marsRows) {
    $c = 1
    foreach ($val in $row) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.Value = "'" + $val
        $cell.ClearFormats()
        $c = $c + 1
    }
    $r = $r + 1
}

# Create the new sheet for April's expenses, right after "mars".
$newWs = $wb.Worksheets.Add($null, $ws)
$newWs.Name = "april"

$newWs.Range("A1").Value = "'Category"
$newWs.Range("B1").Value = "'Name"
$newWs.Range("C1").Value = "'Date"
$newWs.Range("D1").Value = "'Price"
$newWs.Range("E1").Value = "'Account"
$newWs.Range("A1:E1").ClearFormats()

$aprilRow = @("Transportation", "dfsf", "2023-04-12", "342.0", "Checkings")
$c = 1
foreach ($val in $aprilRow) {
    $cell = $newWs.Cells.Item(2, $c)
    $cell.Value = "'" + $val
    $cell.ClearFormats()
    $c = $c + 1
}

# Keep "mars" as the active sheet/tab, matching the original selection.
$ws.Activate()
